$d = $word.ActiveDocument

# -----------------------------------------------------------------------
# Helper: find a unique marker range, then perform a bounded Find/Replace
# immediately after it (avoids Range.Collapse(), which has been observed
# to corrupt the scope of a subsequent Find/Replace-All in this runtime).
# -----------------------------------------------------------------------
function Replace-AfterMarker($markerText, $searchText, $replaceText, $windowSize) {
    $marker = $d.Content
    $found = $marker.Find.Execute($markerText)
    if (-not $found) {
        throw "Marker text not found: $markerText"
    }
    $afterMarker = $marker.End
    $scoped = $d.Range($afterMarker, $afterMarker + $windowSize)
    $ok = $scoped.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 0, $false, $replaceText, 2)
    if (-not $ok) {
        throw "Search text not found after marker '$markerText': $searchText"
    }
}

# -------------------------------------------------------------------
# 1) The three inline maths that just render the number-of-players
#    variable "n" become "N". These live inside <m:oMath> runs, which
#    are opaque to normal text Find/Replace, so we go through the
#    Word OMath object model instead. They are OMaths #1, #8 and #22
#    (each is the lone "n" immediately following "In an ... player
#    normal form game").
# -------------------------------------------------------------------
foreach ($idx in 1, 8, 22) {
    $m = $d.OMaths.Item($idx)
    if ($m.Range.Text -ne "N") {
        $m.Range.Text = "N"
    }
}

# -------------------------------------------------------------------
# 2) "a fixed strategy profile for all other players in the game." ->
#    "an incomplete strategy profile for all other players in the
#    game." with "incomplete strategy" bolded. This must run before
#    the heading-renaming step below, because that step introduces
#    its own (unrelated) occurrence of the phrase "incomplete
#    strategy" in the "Definition of an incomplete strategy profile"
#    heading, which would otherwise be matched instead.
# -------------------------------------------------------------------
$d.Content.Find.Execute("a fixed strategy profile for all other players in the game.", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "an incomplete strategy profile for all other players in the game.", 2) | Out-Null

$boldRange = $d.Content
$boldRange.Find.Execute("incomplete strategy") | Out-Null
$boldRange.Bold = 1

# -------------------------------------------------------------------
# 3) "if there is a strategy (pure or mixed)" -> "if there is a
#    strategy" (both occurrences need the identical change, so a
#    single global replace-all is safe here).
# -------------------------------------------------------------------
$d.Content.Find.Execute("if there is a strategy (pure or mixed)", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "if there is a strategy", 2) | Out-Null

# -------------------------------------------------------------------
# 4) "Definition" -> full heading names (3 occurrences, each scoped to
#    the unique sentence that precedes it so the other two are left
#    untouched).
# -------------------------------------------------------------------
Replace-AfterMarker "To formalise this we need a couple of definitions." `
    "Definition" "Definition of an incomplete strategy profile" 60

Replace-AfterMarker "This notation now allows us to define an important notion in game theory." `
    "Definition" "Definition of a strictly dominated strategy" 60

Replace-AfterMarker "giving the following predicted strategy profile:" `
    "Definition" "Definition of a weakly dominated strategy" 60

# -------------------------------------------------------------------
# 5) Rename the bookmarks that wrap each "Definition" heading so their
#    ids/names match the new heading text.
# -------------------------------------------------------------------
function Rename-Bookmark($oldName, $newName, $markerText) {
    $marker = $d.Content
    $found = $marker.Find.Execute($markerText)
    if ($found) {
        $r = $d.Range($marker.Start, $marker.End)
        if ($d.Bookmarks.Exists($oldName)) {
            $d.Bookmarks.Item($oldName).Delete()
        }
        $d.Bookmarks.Add($newName, $r)
    }
}

Rename-Bookmark "definition" "definition-of-an-incomplete-strategy-profile" "Definition of an incomplete strategy profile"
Rename-Bookmark "definition-1" "definition-of-a-strictly-dominated-strategy" "Definition of a strictly dominated strategy"
Rename-Bookmark "definition-2" "definition-of-a-weakly-dominated-strategy" "Definition of a weakly dominated strategy"
